# B6-PowerPoint.pptx edit — Mon, Jun 29, 2020 10:07:25 AM
#
# 1) The three tables (slides 14, 15, 16) get their custom table style
#    swapped out for a built-in PowerPoint table style.
# 2) The presentation's theme colour scheme ("Integral" / "Red Violet")
#    is replaced with the plain default Office colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------------
$newTableStyleId = "{F32B426D-8A18-4AE9-9B29-0F71E719CDE8}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the deck's theme colours back to the default Office scheme ----
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink — in that order, as VBA
#  OLE-color / BGR-packed RGB() integers.)
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le 12; $k++) {
    $themeColors.Colors($k).RGB = $officeColors[$k - 1]
}
